# bug002_1_1 update:
#  - A12: "在找回密码界面" -> "找回密码界面"
#  - B10: "2020.08.11" -> "2020.08.12"
#  - leave final selection on B13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = "找回密码界面"
$ws.Range("B10").Value = "2020.08.12"

$ws.Range("B13").Select()
